# Actividad-escritura-palabra.xlsx — add translated-phrase column data
# Adds rows for levels 4-10 (song "Smart money..." / "!chan chan chaaann!")
# each paired with a new "Traduccion N" string, fixes the C1 header typo
# ("Fase traducida" -> "Frase traducida"), and leaves the selection on C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header typo on C1 ---
$ws.Range("C1").Value = "Frase traducida"

# --- Song phrase text used on column B (same two strings already used in rows 6-7) ---
$songA = "Smart money bettin' I'll be better off without you"
$songB = "!chan chan chaaann!"

# --- New data rows 8-14 (levels 4-10), alternating the two song-phrase rows,
#     each with the next sequential "Traduccion N" translation ---
$newRows = @(
    @{ Row = 8;  Level = 4;  Phrase = $songA; Style = "A"; Translation = "Traduccion 7" },
    @{ Row = 9;  Level = 5;  Phrase = $songB; Style = "B"; Translation = "Traduccion 8" },
    @{ Row = 10; Level = 6;  Phrase = $songA; Style = "A"; Translation = "Traduccion 9" },
    @{ Row = 11; Level = 7;  Phrase = $songB; Style = "B"; Translation = "Traduccion 10" },
    @{ Row = 12; Level = 8;  Phrase = $songA; Style = "A"; Translation = "Traduccion 11" },
    @{ Row = 13; Level = 9;  Phrase = $songB; Style = "B"; Translation = "Traduccion 12" },
    @{ Row = 14; Level = 10; Phrase = $songA; Style = "A"; Translation = "Traduccion 13" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellB = $ws.Cells.Item($rowNum, 2)
    $cellC = $ws.Cells.Item($rowNum, 3)

    $cellA.Value = $r.Level
    $cellB.Value = $r.Phrase
    $cellC.Value = $r.Translation

    if ($r.Style -eq "A") {
        # Matches the formatting family used on row 6 (A6/B6)
        $cellA.Font.Size = 11
        $cellA.Font.Name = "Calibri"
        $cellA.Font.Bold = $false
        $cellB.Font.Size = 11
        $cellB.Font.Name = "Calibri"
        $cellB.Font.Bold = $false
    } else {
        # Matches the formatting family used on row 7 (A7/B7)
        $cellB.Font.Size = 12
        $cellB.Font.Name = "Calibri"
        $cellB.Font.Bold = $false
        $cellB.Interior.ThemeColor = 3
        $cellB.Interior.TintAndShade = 0.39997558519241921
    }

    # Translation column keeps the same formatting family throughout (row 2-14)
    $cellC.Font.Size = 12
    $cellC.Font.Name = "Calibri"
    $cellC.Font.Bold = $false

    $ws.Rows.Item($rowNum).RowHeight = 15.75
}

# --- Selection ends on C1, matching the authored workbook state ---
$ws.Range("C1").Select() | Out-Null
